# Minor PPT title/text changes
$p = $ppt.ActivePresentation

# Slide 17: "Text 1" shape - remove the "(Example results - adjust based on
# your actual numbers)" italic note entirely (leave the paragraph empty,
# i.e. drop the run rather than just blanking its text).
$s17 = $p.Slides.Item(17)
$tr17 = $s17.Shapes.Item(2).TextFrame.TextRange
$tr17.Characters(1, $tr17.Length).Delete()

# Slide 18: same italic note shape/text as slide 17.
$s18 = $p.Slides.Item(18)
$tr18 = $s18.Shapes.Item(2).TextFrame.TextRange
$tr18.Characters(1, $tr18.Length).Delete()

# Slide 19: update the protocol comparison headline text.
$s19 = $p.Slides.Item(19)
$s19.Shapes.Item(3).TextFrame.TextRange.Text = "Protocol 2 shows ~10-20% higher accuracy than Protocol 1"

# Slide 19: shorten the "Best Practice" guidance text.
$s19.Shapes.Item(16).TextFrame.TextRange.Text = "For UrbanSound8K: Always use predefined folds."

# Slide 21: tweak the summary bullet wording.
$s21 = $p.Slides.Item(21)
$s21.Shapes.Item(11).TextFrame.TextRange.Text = "XGBoost + features tweak: effective, efficient, reproducible approach"
